$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Person" class section (rows 4-8, column D) becomes the "Karyawan" class ---
# D4: Person -> Karyawan (style unchanged, s=10)
$ws.Range("D4").Value = "Karyawan"

# D5: nama -> namaKaryawan (style unchanged, s=2)
$ws.Range("D5").Value = "namaKaryawan"

# D6: alamat -> addKaryawan(), and takes on D7's current highlighted style (s=6)
$ws.Range("D7").Copy()
$ws.Range("D6").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D6").Value = "addKaryawan()"

# D7: addPerson() -> showKaryawan() (style unchanged, s=6)
$ws.Range("D7").Value = "showKaryawan()"

# D8: showPerson() -> empty, takes on the "blank" style used by D17/D18/D19 (s=7)
$ws.Range("D8").Clear()
$ws.Range("D17").Copy()
$ws.Range("D8").PasteSpecial(-4122)   # xlPasteFormats

# --- old "Karyawan" block (rows 13-16, column D) is removed entirely ---
$ws.Range("D13:D16").Clear()

# --- new blank styled cell at I16 (matches styling used by I5:I13) ---
$ws.Range("I5").Copy()
$ws.Range("I16").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I16").ClearContents()

# --- selection moves from E9 to D6 ---
$ws.Range("D6").Select()
